$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "C2" = 29.30314445495605
    "D2" = 0.3231444549560507
    "E2" = 0.1044223387688431
    "B3" = 29.15000000000001
    "C3" = 29.25433921813965
    "D3" = 0.1043392181396428
    "E3" = 0.01088667244199196
    "B4" = 29.34999999999999
    "C4" = 29.70592308044434
    "D4" = 0.3559230804443416
    "E4" = 0.1266812391929893
    "C5" = 29.52533531188965
    "D5" = 0.1553353118896439
    "E5" = 0.02412905911985294
    "B6" = 29.53999999999999
    "C6" = 28.87110710144043
    "D6" = -0.6688928985595624
    "E6" = 0.447417709743413
    "C7" = 29.35497856140137
    "D7" = -0.19502143859863
    "E7" = 0.0380333615130792
    "C8" = 29.55397605895996
    "D8" = -0.1960239410400391
    "E8" = 0.03842538546086871
    "C9" = 30.04802513122559
    "D9" = 0.2080251312255825
    "E9" = 0.04327445522142083
    "C10" = 29.96916007995605
    "D10" = 0.1591600799560524
    "E10" = 0.025331931051617
    "C11" = 29.92043113708496
    "D11" = 0.0004311370849592322
    "E11" = 0.0000001858791860271442
    "C12" = 29.84786033630371
    "D12" = -0.132139663696293
    "E12" = 0.01746089072176943
    "B13" = 30.03999999999999
    "C13" = 30.05769920349121
    "D13" = 0.0176992034912189
    "E13" = 0.0003132618042235752
    "B14" = 30.21000000000001
    "C14" = 30.05219078063965
    "D14" = -0.1578092193603595
    "E14" = 0.02490374971512607
    "C15" = 30.22921562194824
    "D15" = 0.009215621948243324
    "E15" = 0.00008492768789294409
    "C16" = 30.29749870300293
    "D16" = -0.08250129699706577
    "E16" = 0.006806464006198052
    "C17" = 30.6120433807373
    "D17" = 0.172043380737307
    "E17" = 0.02959892485552196
    "C18" = 30.40457725524902
    "D18" = -0.07542274475098054
    "E18" = 0.005688590425771562
    "C19" = 30.41995811462402
    "D19" = -0.2700418853759743
    "E19" = 0.07292261985741083
    "C20" = 30.54005241394043
    "D20" = -0.2099475860595703
    "E20" = 0.04407798889224068
    "C21" = 30.62829399108887
    "D21" = -0.3117060089111305
    "E21" = 0.0971606359913058
    "C22" = 30.73287773132324
    "D22" = -0.2171222686767607
    "E22" = 0.04714207955534344
    "C23" = 31.12346458435059
    "D23" = 0.1034645843505899
    "E23" = 0.01070492021484034
    "C24" = 31.30785179138184
    "D24" = 0.1878517913818314
    "E24" = 0.0352882955253631
    "C25" = 31.35338401794434
    "D25" = 0.0733840179443348
    "E25" = 0.005385214089654452
    "C26" = 31.22195243835449
    "D26" = -0.1580475616455033
    "E26" = 0.02497903174208916
    "C27" = 31.44917106628418
    "D27" = -0.1308289337158186
    "E27" = 0.01711620989721806
    "B28" = 31.65000000000001
    "C28" = 31.89203453063965
    "D28" = 0.2420345306396428
    "E28" = 0.05858071402195217
    "C29" = 32.54359817504883
    "D29" = 0.6635981750488327
    "E29" = 0.4403625379281412
    "C30" = 32.46515274047852
    "D30" = 0.1851527404785145
    "E30" = 0.03428153730670414
    "C31" = 32.5944709777832
    "D31" = 0.1444709777832003
    "E31" = 0.02087186342163395
    "B32" = 32.84999999999999
    "C32" = 32.88671493530273
    "D32" = 0.03671493530274006
    "E32" = 0.001347986474284388
    "B33" = 32.90000000000001
    "C33" = 33.03625106811523
    "D33" = 0.1362510681152287
    "E33" = 0.01856435356254069
    "B34" = 33.09999999999999
    "C34" = 32.95847702026367
    "D34" = -0.1415229797363224
    "E34" = 0.02002875379344753
    "B35" = 33.40000000000001
    "C35" = 33.7332878112793
    "D35" = 0.3332878112792912
    "E35" = 0.1110807651473404
    "C36" = 33.6826286315918
    "D36" = -0.01737136840820597
    "E36" = 0.0003017644403736163
    "B37" = 34.09999999999999
    "C37" = 33.8400993347168
    "D37" = -0.2599006652831974
    "E37" = 0.06754835581464863
    "B38" = 34.40000000000001
    "C38" = 34.36728286743164
    "D38" = -0.03271713256836506
    "E38" = 0.001070410763495974
    "B39" = 34.90000000000001
    "C39" = 34.92761611938477
    "D39" = 0.02761611938475994
    "E39" = 0.0007626500498733137
    "C40" = 35.70381546020508
    "D40" = 0.403815460205081
    "E40" = 0.1630669259006413
    "C41" = 35.95959854125977
    "D41" = 0.2595985412597628
    "E41" = 0.06739140262419677
    "C42" = 36.04191970825195
    "D42" = -0.258080291748044
    "E42" = 0.06660543698875553
    "C43" = 36.59153366088867
    "D43" = -0.2084663391113253
    "E43" = 0.04345821454247807
    "C44" = 37.04032516479492
    "D44" = -0.2596748352050753
    "E44" = 0.067431020038783
    "B45" = 37.90000000000001
    "C45" = 37.8508186340332
    "D45" = -0.04918136596680256
    "E45" = 0.002418806758360565
    "C46" = 38.36010360717773
    "D46" = -0.1398963928222656
    "E46" = 0.01957100072468165
    "B47" = 38.90000000000001
    "C47" = 38.94926071166992
    "D47" = 0.04926071166991619
    "E47" = 0.002426617714226617
    "B48" = 39.40000000000001
    "C48" = 39.48398208618164
    "D48" = 0.08398208618163494
    "E48" = 0.007052990799419559
    "B49" = 39.90000000000001
    "C49" = 39.5393180847168
    "D49" = -0.3606819152832088
    "E49" = 0.1300914440123638
    "B50" = 40.09999999999999
    "C50" = 40.14670181274414
    "D50" = 0.04670181274414631
    "E50" = 0.002181059313589307
    "B51" = 40.59999999999999
    "C51" = 40.69865036010742
    "D51" = 0.09865036010742756
    "E51" = 0.009731893549325134
    "C52" = 0.04815361022947684
    "E52" = 2.654464649066489
    "E53" = 0.05308929298132978
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value2 = $values[$cell]
}
